$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header J1: technical_skill -> skill_group
$ws.Range("J1").Value = "skill_group"

# Replace row 2 with the (former) row 3 candidate's data, with some field tweaks,
# then delete the old row 3 entirely.
$ws.Range("B2").Value = "GOO YE JUI"
$ws.Range("C2").NumberFormat = "@"
$ws.Range("C2").Value = "+60184040438"
$ws.Range("C2").Style = "Normal"
$ws.Range("D2").Value = "yjyejui626@gmail.com"
$ws.Range("E2").Value = "N/A"
$ws.Range("F2").Value = "N/A"
$ws.Range("G2").Value = "[{'Country': 'Malaysia', 'State': 'Penang', 'City': 'Bukit Mertajam'}]"
$ws.Range("H2").Value = "[{'field_of_study': 'Bachelor Of Computer Science (Data Engineering)', 'level': ""Bachelor's Degree"", 'cgpa': '3.97', 'university': 'Universiti Teknologi Malaysia', 'start_date': '2020', 'year_of_graduation': '2024'}, {'field_of_study': 'Foundation in Science', 'level': 'Foundation', 'cgpa': '3.78', 'university': 'Universiti Teknologi Malaysia', 'start_date': '2019', 'year_of_graduation': '2020'}]"
$ws.Range("I2").Value = "['Microsoft Certified: Azure AI Fundamentals', 'Google Data Analytics Certificate by Coursera', 'Alteryx Foundational Micro-Credential', 'Alteryx Designer Core Certification', 'AWS Academy Graduate - AWS Academy Cloud Foundations', 'AWS Academy Graduate - AWS Academy Machine Learning Foundations', 'AWS Academy Graduate - AWS Academy Data Analytics', 'AWS Academy Graduate - AWS Academy Machine Learning for Natural Language Processing', 'AWS Academy Graduate - AWS Academy Data Engineering', 'AWS Academy Graduate - AWS Academy Cloud Web Application Builder', 'AWS Academy Graduate - AWS Academy Cloud Data Pipeline Builder']"
$ws.Range("J2").Value = "['Full-stack web development', 'Natural Language Processing', 'Generative AI']"
$ws.Range("K2").Value = "['HTML 5', 'CSS', 'JavaScript', 'PHP', 'SQL', 'Python', '.NET', 'React', 'spaCy', 'NLTK', 'TensorFlow', 'PyTorch', 'LangChain', 'Llama', 'Django', 'PostgreSQL', 'OpenAI GPT', 'Laravel', 'MySQL', 'Microsoft SQL Server', '.NET MVC Framework']"
$ws.Range("L2").Value = "['English', 'Mandarin', 'Malay', 'French']"
$ws.Range("M2").Value = "[{'job_title': 'Data Science Intern', 'job_company': 'Petronas Digital Sdn Bhd', 'Industries': 'Information Technology', 'start_date': '2023-09', 'end_date': '2024-06', 'job_location': 'N/A'}]"

# Delete row 3 (shift cells up), removing it from the sheet entirely
$ws.Rows.Item(3).Delete()
